$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the date the record was last changed.
# Every data row (2 through 72) is being bumped forward by one day:
# serial 45179 (2023-09-10) -> serial 45180 (2023-09-11).
for ($row = 2; $row -le 72; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45179) {
        $cell.Value2 = 45180
    }
}
